# 01项目计划表.xlsx - "Add files via upload" edit
#
# Summary of the change (per the source diff):
#  - Section "2018.10.17 第七周周三" (E-R 图 建模) 完成情况 column:
#      * 李光洪 (row 53): 未完成 -> 已完成
#      * 吴彤林 (row 54): 未完成 -> 0.8 (80%, entered as a percentage number)
#      * 劳汉文 (row 55): 未完成 -> 已完成
#      * 方嘉耀 (row 56): 未完成 -> 已完成
#      * 成世靖 (row 57): 未完成 -> 已完成
#      * 丰浩   (row 58): 未完成 -> 已完成
#  - Row 59 "总结：" gains the actual summary text.
#  - The "日期" merged-header cells in column A (rows 21/31/41/51) get their
#    font re-resolved (bold, still centered) so they match the style already
#    used by the row 1 / row 11 headers.
#  - The view's active cell moves to D43.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- 完成情况 (column C) updates for the PowerDesigner E-R diagram plan ---
$ws.Range("C53").Value = "已完成"

$ws.Range("C54").Value = 0.8
$ws.Range("C54").NumberFormat = "0%"

$ws.Range("C55").Value = "已完成"
$ws.Range("C56").Value = "已完成"
$ws.Range("C57").Value = "已完成"
$ws.Range("C58").Value = "已完成"

# --- Fill in the final "总结：" (summary) note for this section ---
$ws.Range("A59").Value = "总结：逻辑模型的关系理清、完全完成后才可以进行写物理模型；由于时间关系，物理模型还未整理关系"

# --- Re-resolve the bold/centered font on the other date-header rows so ---
# --- they share the same style record as the first two headers          ---
$ws.Range("A21").Font.Bold = $true
$ws.Range("A31").Font.Bold = $true
$ws.Range("A41").Font.Bold = $true
$ws.Range("A51").Font.Bold = $true

# --- Update the saved view/selection state ---
$excel.ActiveWindow.ScrollRow = 31
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D43").Select()
